$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Prepare & render" re-run: refreshed mean / CI_low / CI_high estimates,
# plus a new grouping column "along" whose value is constant for every row
# in this (extended) variant of the comparison.

$header = @("y", "mean", "CI_low", "CI_high", "along")
for ($c = 0; $c -lt $header.Length; $c++) {
    $ws.Cells.Item(1, $c + 1).Value = $header[$c]
}

$flagValue = "millionaire_tax_in_programTRUE"

# row label, mean, CI_low, CI_high
$rows = @(
    @("<b>All</b>",      0.049590635714958,    0.0352059596767695,  0.0639753117531465),
    @("<b>Europe</b>",   0.0527475039113948,   0.0316466646283374,  0.0738483431944522),
    @("France",          0.0233244079993287,  -0.0286578651114438,  0.0753066811101012),
    @("Germany",         0.0673435539952316,   0.0202569571332121,  0.114430150857251),
    @("Italy",           0.0888154785299581,   0.035649042829676,   0.14198191423024),
    @("Poland",          0.0262053761911532,  -0.0368046718109262,  0.0892154241932327),
    @("Spain",           0.099199542595768,    0.0364259825969358,  0.1619731025946),
    @("United Kingdom",  0.0558479400636768,   0.00221841562454934, 0.109477464502804),
    @("Switzerland",    -0.0172352523098804,  -0.0839669887723522,  0.0494964841525915),
    @("Japan",           0.0305925768083384,   0.0000356375296625966,0.0611495160870142),
    @("USA",             0.0585284294228197,   0.0328795309440684,  0.0841773279015711)
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $rows[$i][0]
    $ws.Cells.Item($r, 2).Value = $rows[$i][1]
    $ws.Cells.Item($r, 3).Value = $rows[$i][2]
    $ws.Cells.Item($r, 4).Value = $rows[$i][3]
    $ws.Cells.Item($r, 5).Value = $flagValue
}
